$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove obsolete rows 23 and 24 (oldest forecast vintage no longer needed) ---
$ws.Rows.Item(23).Delete()
$ws.Rows.Item(23).Delete()

# --- Remove obsolete column BA (oldest forecast horizon column dropped) ---
$ws.Columns.Item(53).Delete()

# --- Update recalculated forecast values ---
$ws.Range("B1").Value = 39583
$ws.Range("C1").Value = 39765
$ws.Range("D1").Value = 39948
$ws.Range("E1").Value = 40130
$ws.Range("F1").Value = 40310
$ws.Range("G1").Value = 40494
$ws.Range("H1").Value = 40676
$ws.Range("I1").Value = 40862
$ws.Range("J1").Value = 41044
$ws.Range("K1").Value = 41228
$ws.Range("L1").Value = 41409
$ws.Range("M1").Value = 41592
$ws.Range("N1").Value = 41774
$ws.Range("O1").Value = 41957
$ws.Range("P1").Value = 42137
$ws.Range("Q1").Value = 42321
$ws.Range("R1").Value = 42503
$ws.Range("S1").Value = 42689
$ws.Range("T1").Value = 42867
$ws.Range("U1").Value = 43053
$ws.Range("V1").Value = 43145
$ws.Range("W1").Value = 43235
$ws.Range("X1").Value = 43326
$ws.Range("Y1").Value = 43418
$ws.Range("Z1").Value = 43510
$ws.Range("AA1").Value = 43600
$ws.Range("AB1").Value = 43691
$ws.Range("AC1").Value = 43783
$ws.Range("AD1").Value = 43875
$ws.Range("AE1").Value = 43966
$ws.Range("AF1").Value = 44068
$ws.Range("AG1").Value = 44159
$ws.Range("AH1").Value = 44251
$ws.Range("AI1").Value = 44341
$ws.Range("AJ1").Value = 44432
$ws.Range("AK1").Value = 44525
$ws.Range("AL1").Value = 44617
$ws.Range("AM1").Value = 44706
$ws.Range("AN1").Value = 44798
$ws.Range("AO1").Value = 44890
$ws.Range("AP1").Value = 44981
$ws.Range("AQ1").Value = 45071
$ws.Range("AR1").Value = 45163
$ws.Range("AS1").Value = 45254
$ws.Range("AT1").Value = 45345
$ws.Range("AU1").Value = 45436
$ws.Range("AV1").Value = 45534
$ws.Range("AW1").Value = 45618
$ws.Range("AX1").Value = 45713
$ws.Range("AY1").Value = 45800
$ws.Range("AZ1").Value = 45891
$ws.Range("B3").Value = 1.516248937663556
$ws.Range("C3").Value = 1.560682679516057
$ws.Range("D3").Value = 1.514319819128396
$ws.Range("B4").Value = 1.507861501681118
$ws.Range("C4").Value = 1.519505601659854
$ws.Range("D4").Value = 1.602279001294704
$ws.Range("E4").Value = 2.1453644888767
$ws.Range("F4").Value = 2.007652128026982
$ws.Range("D5").Value = 1.533634203309853
$ws.Range("E5").Value = 1.862695956043758
$ws.Range("F5").Value = 1.768040115052738
$ws.Range("G5").Value = 1.317672174811868
$ws.Range("H5").Value = 1.022680528298392
$ws.Range("F6").Value = 1.798730406367843
$ws.Range("G6").Value = 1.624644329511749
$ws.Range("H6").Value = 1.510468690286459
$ws.Range("I6").Value = 1.501816644427989
$ws.Range("J6").Value = 0.9070039918702477
$ws.Range("H7").Value = 1.524434521513318
$ws.Range("I7").Value = 1.532369786811083
$ws.Range("J7").Value = 1.042579621507111
$ws.Range("K7").Value = 1.028888107831327
$ws.Range("L7").Value = 1.022042907336096
$ws.Range("J8").Value = 1.052119244508254
$ws.Range("K8").Value = 0.9739397072873635
$ws.Range("L8").Value = 1.083482333436536
$ws.Range("M8").Value = 1.303605130836716
$ws.Range("N8").Value = 1.34489417553354
$ws.Range("L9").Value = 1.048230161591079
$ws.Range("M9").Value = 1.159709768734829
$ws.Range("N9").Value = 1.095916825800991
$ws.Range("O9").Value = 1.192378712846454
$ws.Range("P9").Value = 1.277042522796856
$ws.Range("N10").Value = 1.09266761031781
$ws.Range("O10").Value = 1.224048560391644
$ws.Range("P10").Value = 1.293136192195643
$ws.Range("Q10").Value = 1.210961441871872
$ws.Range("R10").Value = 2.615369162917314
$ws.Range("P11").Value = 1.276496651730441
$ws.Range("Q11").Value = 1.268509999927447
$ws.Range("R11").Value = 1.825134644920934
$ws.Range("S11").Value = 2.033218171624651
$ws.Range("T11").Value = 2.618053282882693
$ws.Range("R12").Value = 1.764874365638147
$ws.Range("S12").Value = 1.931144489665093
$ws.Range("T12").Value = 2.406099663413808
$ws.Range("U12").Value = 2.152537330144288
$ws.Range("V12").Value = 1.731723847815725
$ws.Range("W12").Value = 1.431088640641853
$ws.Range("X12").Value = 1.372961566907027
$ws.Range("T13").Value = 2.338011880246493
$ws.Range("U13").Value = 2.357106719845503
$ws.Range("V13").Value = 2.356276715023498
$ws.Range("W13").Value = 2.21629047761287
$ws.Range("X13").Value = 1.902399534782662
$ws.Range("Y13").Value = 1.966855307908655
$ws.Range("Z13").Value = 1.983559881711905
$ws.Range("AA13").Value = 2.136062314641141
$ws.Range("AB13").Value = 2.241561867365394
$ws.Range("W14").Value = 2.310526014185643
$ws.Range("X14").Value = 1.969642719257858
$ws.Range("Y14").Value = 1.937894966683062
$ws.Range("Z14").Value = 1.93103453922987
$ws.Range("AA14").Value = 2.031764787322499
$ws.Range("AB14").Value = 2.135927826705641
$ws.Range("AC14").Value = 1.950353221540246
$ws.Range("AD14").Value = 2.041276490941102
$ws.Range("AE14").Value = 2.166968775134936
$ws.Range("AF14").Value = 2.139672475020404
$ws.Range("AA15").Value = 1.98865454260444
$ws.Range("AB15").Value = 1.988654542604462
$ws.Range("AC15").Value = 1.806095414188391
$ws.Range("AD15").Value = 1.854752869950294
$ws.Range("AE15").Value = 1.984987808509886
$ws.Range("AF15").Value = 2.011395609719546
$ws.Range("AG15").Value = 2.210985773414453
$ws.Range("AH15").Value = 2.330672672271739
$ws.Range("AI15").Value = 2.459440348120401
$ws.Range("AJ15").Value = 2.100991693542231
$ws.Range("AE16").Value = 1.95617858661592
$ws.Range("AF16").Value = 2.151090884613516
$ws.Range("AG16").Value = 2.232402359458496
$ws.Range("AH16").Value = 2.307457288603798
$ws.Range("AI16").Value = 2.526389380645511
$ws.Range("AJ16").Value = 1.091147151778871
$ws.Range("AK16").Value = 1.114171399050901
$ws.Range("AL16").Value = 0.7608230790701942
$ws.Range("AM16").Value = 0.8171929556848756
$ws.Range("AN16").Value = 0.8967077601845341
$ws.Range("AH17").Value = 2.273076295481702
$ws.Range("AI17").Value = 2.395170926915613
$ws.Range("AJ17").Value = 1.800020576268269
$ws.Range("AK17").Value = 1.718967187759834
$ws.Range("AL17").Value = 1.459415358104388
$ws.Range("AM17").Value = 1.509741350988136
$ws.Range("AN17").Value = 1.773412413757813
$ws.Range("AO17").Value = 0.1338254721205745
$ws.Range("AP17").Value = -0.4883557973630492
$ws.Range("AQ17").Value = 0.8766015904249524
$ws.Range("AL18").Value = 1.557009346267302
$ws.Range("AM18").Value = 1.548916348216034
$ws.Range("AN18").Value = 1.589386175889063
$ws.Range("AO18").Value = 0.9836073856501315
$ws.Range("AP18").Value = 0.4753196237801127
$ws.Range("AQ18").Value = 2.477445663648559
$ws.Range("AR18").Value = 2.01691766737
$ws.Range("AS18").Value = 1.823564868738359
$ws.Range("AT18").Value = 1.617535832906758
$ws.Range("AU18").Value = 1.554086551645839
$ws.Range("AV18").Value = 1.508385007449875
$ws.Range("AP19").Value = 0.6415376390553895
$ws.Range("AQ19").Value = 1.488234279941625
$ws.Range("AR19").Value = 1.205410808100038
$ws.Range("AS19").Value = 0.9634615329177532
$ws.Range("AT19").Value = 1.142484412546874
$ws.Range("AU19").Value = 0.9888012784191602
$ws.Range("AV19").Value = 0.8780954941978392
$ws.Range("AW19").Value = 0.6266145540918089
$ws.Range("AX19").Value = 0.3854686824285025
$ws.Range("AY19").Value = 0.5837948599211717
$ws.Range("AT20").Value = 1.245777925635272
$ws.Range("AU20").Value = 1.124712786946613
$ws.Range("AV20").Value = 1.165055762672873
$ws.Range("AW20").Value = 1.003756253906252
$ws.Range("AX20").Value = 0.6767639290315763
$ws.Range("AY20").Value = 1.328924132093245
$ws.Range("AZ20").Value = 1.609787824259601
$ws.Range("AX21").Value = 0.7219463338497878
$ws.Range("AY21").Value = 1.326548386659265
$ws.Range("AZ21").Value = 1.690613666316931

# --- Clear cells that are no longer part of the forecast window for their vintage ---
$ws.Range("C5").ClearContents()
$ws.Range("E6").ClearContents()
$ws.Range("G7").ClearContents()
$ws.Range("I8").ClearContents()
$ws.Range("K9").ClearContents()
$ws.Range("M10").ClearContents()
$ws.Range("O11").ClearContents()
$ws.Range("Q12").ClearContents()
$ws.Range("R13").ClearContents()
$ws.Range("S13").ClearContents()
$ws.Range("T14").ClearContents()
$ws.Range("U14").ClearContents()
$ws.Range("V14").ClearContents()
$ws.Range("V15").ClearContents()
$ws.Range("W15").ClearContents()
$ws.Range("X15").ClearContents()
$ws.Range("Y15").ClearContents()
$ws.Range("Z15").ClearContents()
$ws.Range("Y16").ClearContents()
$ws.Range("Z16").ClearContents()
$ws.Range("AA16").ClearContents()
$ws.Range("AB16").ClearContents()
$ws.Range("AC16").ClearContents()
$ws.Range("AD16").ClearContents()
$ws.Range("AC17").ClearContents()
$ws.Range("AD17").ClearContents()
$ws.Range("AE17").ClearContents()
$ws.Range("AF17").ClearContents()
$ws.Range("AG17").ClearContents()
$ws.Range("AG18").ClearContents()
$ws.Range("AH18").ClearContents()
$ws.Range("AI18").ClearContents()
$ws.Range("AJ18").ClearContents()
$ws.Range("AK18").ClearContents()
$ws.Range("AK19").ClearContents()
$ws.Range("AL19").ClearContents()
$ws.Range("AM19").ClearContents()
$ws.Range("AN19").ClearContents()
$ws.Range("AO19").ClearContents()
$ws.Range("AO20").ClearContents()
$ws.Range("AP20").ClearContents()
$ws.Range("AQ20").ClearContents()
$ws.Range("AR20").ClearContents()
$ws.Range("AS20").ClearContents()
$ws.Range("AS21").ClearContents()
$ws.Range("AT21").ClearContents()
$ws.Range("AU21").ClearContents()
$ws.Range("AV21").ClearContents()
$ws.Range("AW21").ClearContents()
$ws.Range("AW22").ClearContents()
$ws.Range("AX22").ClearContents()
$ws.Range("AY22").ClearContents()
$ws.Range("AZ22").ClearContents()
